# Update "想去人数" (F column) counts across all sheets to reflect the
# latest scrape output (gh-pages regeneration at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 899
$ws1.Range("F4").Value = 4354
$ws1.Range("F7").Value = 3458
$ws1.Range("F8").Value = 982
$ws1.Range("F12").Value = 2351
$ws1.Range("F13").Value = 1253
$ws1.Range("F16").Value = 506
$ws1.Range("F19").Value = 9740
$ws1.Range("F20").Value = 5928
$ws1.Range("F23").Value = 812
$ws1.Range("F26").Value = 3513
$ws1.Range("F34").Value = 17
$ws1.Range("F37").Value = 19

# 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F22").Value = 74

# 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8646
$ws3.Range("F4").Value = 1519

# 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8646
$ws4.Range("F3").Value = 899
$ws4.Range("F5").Value = 1519
$ws4.Range("F7").Value = 4354
$ws4.Range("F10").Value = 3458
$ws4.Range("F11").Value = 982
$ws4.Range("F15").Value = 2351
$ws4.Range("F19").Value = 1253
$ws4.Range("F23").Value = 506
$ws4.Range("F26").Value = 9740
$ws4.Range("F31").Value = 812
$ws4.Range("F34").Value = 3513
$ws4.Range("F46").Value = 74
